$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered - same as existing headers) into the new header cells
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

# Set header text
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I and J values per data row (2-21)
$IVals = @{2=1;3=1;4=1;5=1;6=1;7=1;8=1;9=1;10=1;11=1;12=5;13=1;14=1;15=1;16=1;17=1;18=1;19=1;20=1;21=2}
$JVals = @{2=4;3=7;4=4;5=6;6=7;7=6;8=7;9=6;10=5;11=4;12=6;13=7;14=8;15=5;16=5;17=5;18=4;19=3;20=2;21=3}

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 9).Value = $IVals[$r]
    $ws.Cells.Item($r, 10).Value = $JVals[$r]
}
